$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22; everything from row 22 down shifts to row 23+
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new observation
$ws.Cells.Item(22, 1).Value = 10
$ws.Cells.Item(22, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(22, 3).Value = "La Araucanía"
$ws.Cells.Item(22, 4).Value = 44623
$ws.Cells.Item(22, 5).Value = 9
$ws.Cells.Item(22, 6).Value = 100112030
$ws.Cells.Item(22, 7).Value = "Poroto granado"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 90
$ws.Cells.Item(22, 11).Value = 25000
$ws.Cells.Item(22, 12).Value = 25000
$ws.Cells.Item(22, 13).Value = 25000
$ws.Cells.Item(22, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(22, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(22, 16).Value = 1000
$ws.Cells.Item(22, 17).Value = 25
$ws.Cells.Item(22, 18).Value = "Hortaliza"
